# "copy changes to template" -- clear the per-field column headers (keep
# only the "ID" label), drop the bold weight from the banner rows, add a
# blank spacer row above the header row, set the print/page setup, and
# move the active selection down below the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Header row (row 7): keep "ID" in A7 but blank out the rest of the
#    column captions (Lastname, Firstname, Middle, Suffix, Birthdate,
#    Marital, Sex, Religion, Nationality, Occupation, Status, Date record).
#    Clearing the contents also prunes the now-unused shared strings.
$ws.Range("B7:M7").ClearContents()

# 2. The banner/title rows (1-4, each merged A:M) lose their bold weight.
$ws.Range("A1:M4").Font.Bold = $false

# 3. Insert a blank formatted spacer row at row 6 (bold, centered) between
#    the banner block and the header row.
$row6 = $ws.Rows.Item(6)
$row6.Font.Bold = $true
$row6.HorizontalAlignment = -4108
$row6.VerticalAlignment = -4108
$row6.WrapText = $false

# 4. Page setup: A4 paper, portrait orientation.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# 5. Move the active selection to A8 (just below the header row).
$ws.Range("A8").Select() | Out-Null
